$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.02110846634833865
$ws.Range("E4").Value = 0.01298645735145829
$ws.Range("F4").Value = 0.01845237328209493
$ws.Range("H4").Value = -0.01932192278887691
$ws.Range("J4").Value = 0.00912542055411643
$ws.Range("C5").Value = -0.01038802927952117
$ws.Range("E5").Value = 0.0008490873939634957
$ws.Range("F5").Value = -0.0008634068505362739
$ws.Range("H5").Value = 0.001045050665802027
$ws.Range("J5").Value = -0.005236802687123067
$ws.Range("C6").Value = 0.005607882560315302
$ws.Range("E6").Value = -0.01080562968022519
$ws.Range("F6").Value = 0.004349788397991535
$ws.Range("H6").Value = -0.004340060909602436
$ws.Range("J6").Value = 0.004516980638787003
$ws.Range("C7").Value = 0.00371736101269444
$ws.Range("E7").Value = -0.006509300324372012
$ws.Range("F7").Value = -0.005087259371490373
$ws.Range("H7").Value = 0.00526512376260495
$ws.Range("J7").Value = 0.01708004907309691
$ws.Range("C8").Value = 0.1006796693711868
$ws.Range("E8").Value = -0.1291824037272961
$ws.Range("F8").Value = -0.9990109740724388
$ws.Range("H8").Value = 0.9999999321279971
$ws.Range("J8").Value = -0.04652471361905297
$ws.Range("C9").Value = 0.9555876790235071
$ws.Range("E9").Value = 0.009469105242764208
$ws.Range("F9").Value = 0.02704875573795023
$ws.Range("H9").Value = -0.02724844620993785
$ws.Range("J9").Value = -0.001473274547790203
$ws.Range("C10").Value = 0.004201949160077966
$ws.Range("E10").Value = -0.02559751465590058
$ws.Range("F10").Value = -0.01065525652221026
$ws.Range("H10").Value = 0.01121473408058936
$ws.Range("J10").Value = 0.009003673972347252
$ws.Range("C11").Value = -0.002770926926837077
$ws.Range("E11").Value = 0.01180795314431812
$ws.Range("F11").Value = -0.005739596005583839
$ws.Range("H11").Value = 0.005149507405980296
$ws.Range("J11").Value = 0.03298319604652162
$ws.Range("C12").Value = 0.06583467911338714
$ws.Range("E12").Value = 0.005426115865044634
$ws.Range("F12").Value = -0.02376828920673157
$ws.Range("H12").Value = 0.0238435650177426
$ws.Range("J12").Value = 0.005048352938450342
$ws.Range("C13").Value = 0.1039959086718363
$ws.Range("E13").Value = -0.007767489814699592
$ws.Range("F13").Value = 0.01314472449378898
$ws.Range("H13").Value = -0.01319980689599227
$ws.Range("J13").Value = -0.0005701241463470946
$ws.Range("C14").Value = -0.2008365560334622
$ws.Range("E14").Value = -0.006794817487792698
$ws.Range("F14").Value = -0.0204392476335699
$ws.Range("H14").Value = 0.02071611538864461
$ws.Range("J14").Value = 0.02110973848436875
$ws.Range("C15").Value = -0.01358717593548704
$ws.Range("E15").Value = 0.007514182764567309
$ws.Range("F15").Value = -0.003017292696691707
$ws.Range("H15").Value = 0.002034027633361105
$ws.Range("J15").Value = 0.009964063139056177
$ws.Range("C16").Value = -0.005732050405282015
$ws.Range("E16").Value = 0.02061093135243725
$ws.Range("F16").Value = -0.02791832166073286
$ws.Range("H16").Value = 0.02683639825745593
$ws.Range("J16").Value = -0.006952602142632662
$ws.Range("C17").Value = 0.008226549065061961
$ws.Range("E17").Value = -0.01961454337658173
$ws.Range("F17").Value = -0.04102408320896332
$ws.Range("H17").Value = 0.04119682091187282
$ws.Range("J17").Value = -0.005043380102597345
$ws.Range("C18").Value = 0.02770170610006824
$ws.Range("E18").Value = -0.007052324922092996
$ws.Range("F18").Value = -0.01366606067464243
$ws.Range("H18").Value = 0.01337110843884434
$ws.Range("J18").Value = -0.01201817090873828
$ws.Range("C19").Value = 0.0113806524232261
$ws.Range("E19").Value = 0.008526018101040722
$ws.Range("F19").Value = -0.002496572067862882
$ws.Range("H19").Value = 0.002056993330279733
$ws.Range("J19").Value = -0.02763781249191481
$ws.Range("C20").Value = 0.009769741446789656
$ws.Range("E20").Value = 0.02642357990494319
$ws.Range("F20").Value = -0.0001738172229526889
$ws.Range("H20").Value = -0.0005815046632601864
$ws.Range("J20").Value = 0.001269696444360062
$ws.Range("C21").Value = 0.02399470425578816
$ws.Range("E21").Value = -0.02289703637188145
$ws.Range("F21").Value = -0.02384739839389593
$ws.Range("H21").Value = 0.02410911158836446
$ws.Range("J21").Value = -0.02232792274893487
$ws.Range("C22").Value = 0.01143177376927095
$ws.Range("E22").Value = 0.01954322622172905
$ws.Range("F22").Value = 0.0004399732975989318
$ws.Range("H22").Value = -0.0007685764147430565
$ws.Range("J22").Value = -0.009056289951286765
$ws.Range("C23").Value = -0.008521247668849907
$ws.Range("E23").Value = -0.001891189899647596
$ws.Range("F23").Value = 0.009148773005950919
$ws.Range("H23").Value = -0.008496730227869209
$ws.Range("J23").Value = 0.0008275322108500017
